$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on cells whose new values look like plain numbers,
# so Excel keeps them as literal text instead of converting to a numeric value
# (this preserves exact formatting such as trailing zeros).
$textCells = @('D5','D6','D9','D10','D16','D19','D23','D25','D27','D29','D32','D34','D37','D39','D40','D43','D44','D46','D47','D48','D49')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all cell value updates
$ws.Range('D2').Value = '26.242.28'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '1.591.41'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '212.47'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').Value = '0.499'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('D9').Value = '0.0605'
$ws.Range('E9').Value = '  -0.55%  '
$ws.Range('D10').Value = '18.94'
$ws.Range('E10').Value = '  -2.13%  '
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').Value = '1.815.78'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').Value = '1.587.02'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('E14').Value = '  -1.25%  '
$ws.Range('E15').Value = '  -2.65%  '
$ws.Range('D16').Value = '63.84'
$ws.Range('E16').Value = '  -0.96%  '
$ws.Range('D17').Value = '26.250.08'
$ws.Range('D18').Value = '0.0₃0721'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').Value = '214.87'
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('E20').Value = '  -1.57%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').Value = '8.99'
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('E24').Value = '  -3.04%  '
$ws.Range('D25').Value = '145.31'
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').Value = '6.96'
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('D29').Value = '15.10'
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').Value = '3.19'
$ws.Range('E32').Value = '  -0.48%  '
$ws.Range('D33').Value = '1.416.31'
$ws.Range('E33').Value = '  +5.60%  '
$ws.Range('D34').Value = '2.95'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('E35').Value = '  -0.94%  '
$ws.Range('E36').Value = '  -1.39%  '
$ws.Range('D37').Value = '0.567'
$ws.Range('E37').Value = '  -4.87%  '
$ws.Range('E38').Value = '  -0.53%  '
$ws.Range('D39').Value = '0.823'
$ws.Range('E39').Value = '  +0.75%  '
$ws.Range('D40').Value = '5.78'
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  +0.97%  '
$ws.Range('D43').Value = '0.936'
$ws.Range('E43').Value = '  -7.45%  '
$ws.Range('D44').Value = '0.760'
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('D45').Value = '1.728.60'
$ws.Range('D46').Value = '60.79'
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('D47').Value = '86.68'
$ws.Range('E47').Value = '  -0.75%  '
$ws.Range('D48').Value = '1.47'
$ws.Range('E48').Value = '  -1.63%  '
$ws.Range('D49').Value = '0.0500'
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('E50').Value = '  -3.09%  '
$ws.Range('E51').Value = '  -0.02%  '

# Restore default style on the cells we temporarily formatted as text,
# so no stray style attribute is left on the cell itself.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
